$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - C2 and D2 change
$ws.Range("C2").Value = 6353912774558873
$ws.Range("D2").Value = 6353912774558873

# Row 3 (RandomForestRegressor) - B3, C3, D3 change
$ws.Range("B3").Value = 236267315465153.4
$ws.Range("C3").Value = 226316283694391.1
$ws.Range("D3").Value = 868639828996547.1

# Row 4 (GradientBoostingRegressor -> DecisionTreeRegressor)
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 9955873257502.84
$ws.Range("C4").Value = 10060287319479.29
$ws.Range("D4").Value = 270922304131884.2

# Row 5 (AdaBoostRegressor -> MLPRegressor)
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 430491909386283.6
$ws.Range("C5").Value = 559469640814184.7
$ws.Range("D5").Value = 3195131380174678
